# Applies the diff: clears the final paragraph's run text ("OFFLOAD SPAWNING...")
# and appends a sequence of new paragraphs (notes/todos) before the section break.

$d = $word.ActiveDocument

# --- Step 1: clear the text of the last paragraph, keeping its paragraph mark/props ---
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$clearRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$clearRange.Text = ""

# --- Helper-style sequence: each new paragraph is appended after the current last one ---

# 2: empty paragraph
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()

# 3: "Spawn Pool isn't updating correctly..."
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Spawn Pool isn" + [char]8217 + "t updating correctly on secondary pool objects. So when they die, their children don" + [char]8217 + "t go in the tertiary pool."

# 4: "Spawn zone objects are still spawning in overlapping."
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Spawn zone objects are still spawning in overlapping."

# 5: "Refactor spawning - break out spawners."  (en dash)
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Refactor spawning " + [char]8211 + " break out spawners."

# 6: "Children are sometimes spawning in overlapping."
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Children are sometimes spawning in overlapping."

# 7: "Physics looks bad for asteroids. Is there a physics material that would improve things?"
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Physics looks bad for asteroids. Is there a physics material that would improve things?"

# 8: empty paragraph
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()

# 9: "Refactor Actor components."
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Refactor Actor components."

# 10: empty paragraph
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()

# 11: "Add modifiers for asteroid spawn directions/behaviors."
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Add modifiers for asteroid spawn directions/behaviors."

# 12: tab + "One way (which direction)"
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = [char]9 + "One way (which direction)"

# 13: tab + "Orbiting center" + " (clockwise/counter)" (two runs)
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = [char]9 + "Orbiting center"
$tail = $d.Range($cur.Range.End, $cur.Range.End)
$tail.InsertAfter(" (clockwise/counter)")
$tail.Bold = 1
$tail.Bold = 0

# 14: tab + "Pulsing out in"
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = [char]9 + "Pulsing out in"

# 15: empty paragraph
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()

# 16: "Improve piloting control. Impulse."
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Improve piloting control. Impulse."

# 17: empty paragraph
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()

# 18: "Add weapons."
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs.Last
$cur.Range.Text = "Add weapons."

# 19: final empty paragraph (before sectPr)
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
